# Applies commit "#5: insurance, claim, debt, investment done"
# - Fixes the 債務 (Debt) sheet header row (previously duplicated row-2 data)
# - Fixes the 事業投資 (Investment) sheet header row (previously owner/company data)
# - Normalises the E2 amount cell on both sheets from text to a real number
# - Appends the common trailing metadata columns (property_category, category,
#   date, legislator_name, legislator_id, source_file, index) that already exist
#   on the 土地 (Land) sheet, matching columns H:N there.

$wb = $excel.ActiveWorkbook

# ===================== 債務 (Debt) sheet =====================
$ws2 = $wb.Worksheets.Item("債務")

# Row 1 on this sheet was data (copy of row 2), not headers. Replace B1:G1
# with the real field names.
$ws2Header = @('species', 'debtor', 'owner', 'total', 'register_date', 'register_reason')
for ($i = 0; $i -lt $ws2Header.Length; $i++) {
    $ws2.Cells.Item(1, $i + 2).Value = $ws2Header[$i]
}

# E2 ('total') was stored as text; store it as a real number instead.
$ws2.Cells.Item(2, 5).Value = 5707475

# Extend the table with the trailing metadata columns H:N, copying the
# existing header/body formatting so the new cells match s=1 / s=2.
$ws2.Cells.Item(1, 2).Copy() | Out-Null
$ws2.Range($ws2.Cells.Item(1, 8), $ws2.Cells.Item(1, 14)).PasteSpecial(-4122) | Out-Null
$ws2.Cells.Item(2, 2).Copy() | Out-Null
$ws2.Range($ws2.Cells.Item(2, 8), $ws2.Cells.Item(10, 14)).PasteSpecial(-4122) | Out-Null

$ws2NewHeader = @('property_category', 'category', 'date', 'legislator_name', 'legislator_id', 'source_file', 'index')
for ($c = 0; $c -lt $ws2NewHeader.Length; $c++) {
    $ws2.Cells.Item(1, $c + 8).Value = $ws2NewHeader[$c]
}

$ws2NewRows = @(
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 84),
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 85),
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 86),
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 88),
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 89),
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 90),
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 91),
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 92),
    @('debt', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 93)
)
for ($r = 0; $r -lt $ws2NewRows.Length; $r++) {
    $rowVals = $ws2NewRows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 8).Value = $rowVals[$c]
    }
}

# ===================== 事業投資 (Investment) sheet =====================
$ws3 = $wb.Worksheets.Item("事業投資")

# Row 1 held owner/company sample data instead of headers; replace B1:G1.
$ws3Header = @('owner', 'company', 'address', 'total', 'register_date', 'register_reason')
for ($i = 0; $i -lt $ws3Header.Length; $i++) {
    $ws3.Cells.Item(1, $i + 2).Value = $ws3Header[$i]
}

# E2 ('total') was stored as text; store it as a real number instead.
$ws3.Cells.Item(2, 5).Value = 1676000

# Extend the table with the trailing metadata columns H:N, copying the
# existing header/body formatting so the new cells match s=1 / s=2.
$ws3.Cells.Item(1, 2).Copy() | Out-Null
$ws3.Range($ws3.Cells.Item(1, 8), $ws3.Cells.Item(1, 14)).PasteSpecial(-4122) | Out-Null
$ws3.Cells.Item(2, 2).Copy() | Out-Null
$ws3.Range($ws3.Cells.Item(2, 8), $ws3.Cells.Item(4, 14)).PasteSpecial(-4122) | Out-Null

$ws3NewHeader = @('property_category', 'category', 'date', 'legislator_name', 'legislator_id', 'source_file', 'index')
for ($c = 0; $c -lt $ws3NewHeader.Length; $c++) {
    $ws3.Cells.Item(1, $c + 8).Value = $ws3NewHeader[$c]
}

$ws3NewRows = @(
    @('investment', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 98),
    @('investment', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 99),
    @('investment', 'normal', '2012-11-28', '顏清標', 979, 'tmp68961', 100)
)
for ($r = 0; $r -lt $ws3NewRows.Length; $r++) {
    $rowVals = $ws3NewRows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws3.Cells.Item($r + 2, $c + 8).Value = $rowVals[$c]
    }
}

